$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111523657
$ws.Range("B2").Value = 89686
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 658
$ws.Range("F2").Value = 'Rosenticka'
$ws.Range("G2").Value = 'Rhodofomes roseus'
$ws.Range("H2").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q2").Value = 497390.1961838813
$ws.Range("R2").Value = 6754097.842248607

# Row 3
$ws.Range("A3").Value = 111523727
$ws.Range("B3").Value = 89845
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 1209
$ws.Range("F3").Value = 'Rynkskinn'
$ws.Range("G3").Value = 'Phlebia centrifuga'
$ws.Range("H3").Value = 'P.Karst.'
$ws.Range("Q3").Value = 497338.5868253836
$ws.Range("R3").Value = 6754122.194367126

# Row 4
$ws.Range("A4").Value = 111523656
$ws.Range("B4").Value = 89845
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 1209
$ws.Range("F4").Value = 'Rynkskinn'
$ws.Range("G4").Value = 'Phlebia centrifuga'
$ws.Range("H4").Value = 'P.Karst.'
$ws.Range("Q4").Value = 497390.1961838813
$ws.Range("R4").Value = 6754097.842248607

# Row 5
$ws.Range("A5").Value = 111523741
$ws.Range("B5").Value = 89686
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 658
$ws.Range("F5").Value = 'Rosenticka'
$ws.Range("G5").Value = 'Rhodofomes roseus'
$ws.Range("H5").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q5").Value = 497384.3941364431
$ws.Range("R5").Value = 6754155.713205664

# Row 6
$ws.Range("A6").Value = 111523731
$ws.Range("B6").Value = 89845
$ws.Range("D6").Value = 'VU'
$ws.Range("E6").Value = 1209
$ws.Range("F6").Value = 'Rynkskinn'
$ws.Range("G6").Value = 'Phlebia centrifuga'
$ws.Range("H6").Value = 'P.Karst.'
$ws.Range("Q6").Value = 497307.3714758331
$ws.Range("R6").Value = 6754063.864355386

# Row 7
$ws.Range("A7").Value = 111523683
$ws.Range("B7").Value = 89845
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 1209
$ws.Range("F7").Value = 'Rynkskinn'
$ws.Range("G7").Value = 'Phlebia centrifuga'
$ws.Range("H7").Value = 'P.Karst.'
$ws.Range("Q7").Value = 497391.6869587752
$ws.Range("R7").Value = 6754138.20205555

# Row 8
$ws.Range("A8").Value = 111523728
$ws.Range("B8").Value = 89686
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 658
$ws.Range("F8").Value = 'Rosenticka'
$ws.Range("G8").Value = 'Rhodofomes roseus'
$ws.Range("H8").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q8").Value = 497338.5868253836
$ws.Range("R8").Value = 6754122.194367126

# Row 9
$ws.Range("A9").Value = 111523695
$ws.Range("B9").Value = 5113
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 100526
$ws.Range("F9").Value = 'Bronshjon'
$ws.Range("G9").Value = 'Callidium coriaceum'
$ws.Range("H9").Value = 'Paykull, 1800'
$ws.Range("Q9").Value = 497354.1644349985
$ws.Range("R9").Value = 6754111.484663551

# Row 10
$ws.Range("A10").Value = 111523712
$ws.Range("B10").Value = 89405
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = 'Ullticka'
$ws.Range("G10").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H10").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q10").Value = 497301.0581945881
$ws.Range("R10").Value = 6754088.183226726

# Row 11
$ws.Range("A11").Value = 111523730
$ws.Range("B11").Value = 89405
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = 'Ullticka'
$ws.Range("G11").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H11").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q11").Value = 497338.5868253836
$ws.Range("R11").Value = 6754122.194367126

# Row 12
$ws.Range("A12").Value = 111523701
$ws.Range("B12").Value = 89686
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 658
$ws.Range("F12").Value = 'Rosenticka'
$ws.Range("G12").Value = 'Rhodofomes roseus'
$ws.Range("H12").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q12").Value = 497367.2942720717
$ws.Range("R12").Value = 6754083.757028132

# Row 13
$ws.Range("A13").Value = 111523740
$ws.Range("B13").Value = 56398
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = 'Tretåig hackspett'
$ws.Range("G13").Value = 'Picoides tridactylus'
$ws.Range("H13").Value = '(Linnaeus, 1758)'
$ws.Range("Q13").Value = 497366.3615979423
$ws.Range("R13").Value = 6754139.679549156

# Row 14
$ws.Range("A14").Value = 111523724
$ws.Range("B14").Value = 93881
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 2869
$ws.Range("F14").Value = 'Bollvitmossa'
$ws.Range("G14").Value = 'Sphagnum wulfianum'
$ws.Range("H14").Value = 'Girg.'
$ws.Range("Q14").Value = 497291.3182300103
$ws.Range("R14").Value = 6754089.649475355

# Row 15
$ws.Range("A15").Value = 111523697
$ws.Range("B15").Value = 77515
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = 'Garnlav'
$ws.Range("G15").Value = 'Alectoria sarmentosa'
$ws.Range("H15").Value = '(Ach.) Ach.'
$ws.Range("Q15").Value = 497380.5053056676
$ws.Range("R15").Value = 6754165.927741241

# Move the 'Gamla hack' public comment from row 5 to row 13
$ws.Range("AC5").ClearContents()
$ws.Range("AC13").Value = 'Gamla hack'
